$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# New layout: Name | Course | [Type | Missing Prerequisite N | Reason] x12
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Course"

$col = 3
for ($i = 1; $i -le 12; $i++) {
    $ws.Cells.Item(1, $col).Value = "Type"
    $ws.Cells.Item(1, $col + 1).Value = "Missing Prerequisite $i"
    $ws.Cells.Item(1, $col + 2).Value = "Reason"
    $col = $col + 3
}

# Re-apply the bold+bordered header style across the full new header range (A1:AL1)
$ws.Range("A1:AL1").Borders.LineStyle = 1
$ws.Range("A1:AL1").Font.Bold = $true

# --- Data rows ---
# Row 2: EN 110 / EN 103.0 / Missing
$ws.Cells.Item(2, 1).Value = "Elettra Scianetti"
$ws.Cells.Item(2, 2).Value = "EN 110"
$ws.Cells.Item(2, 3).Value = "prerequisite"
$ws.Cells.Item(2, 4).Value = "EN 103.0"
$ws.Cells.Item(2, 5).Value = "Missing"

# Row 3: EN 110 / EN 105.0 / Missing (new row, replaces old blank row 3/4 and shifts old row5 data)
$ws.Cells.Item(3, 1).Value = "Elettra Scianetti"
$ws.Cells.Item(3, 2).Value = "EN 110"
$ws.Cells.Item(3, 3).Value = "prerequisite"
$ws.Cells.Item(3, 4).Value = "EN 105.0"
$ws.Cells.Item(3, 5).Value = "Missing"

# Row 4: FIN 372 / FIN 301.0 / Grade (was row 5, now moved up to row 4)
$ws.Cells.Item(4, 1).Value = "Elettra Scianetti"
$ws.Cells.Item(4, 2).Value = "FIN 372"
$ws.Cells.Item(4, 3).Value = "prerequisite"
$ws.Cells.Item(4, 4).Value = "FIN 301.0"
$ws.Cells.Item(4, 5).Value = "Grade"

# Apply the bordered data-row style to the new data block (A2:E4)
$ws.Range("A2:E4").Borders.LineStyle = 1

# Column E width (auto best-fit sized to its "Missing"/"Grade" content, ~7.71 chars)
$ws.Columns.Item(5).ColumnWidth = 6.75

# Remove the now-unused old row 5 entirely (its data was moved up into row 4 above)
$ws.Rows.Item(5).Delete()
